# The target revision only re-serializes the package's OOXML parts
# (word/document.xml, word/footnotes.xml, word/header1.xml,
# word/styles.xml): every changed line is strictly an attribute /
# namespace-declaration reordering (e.g. w:type="default" r:id="rId6"
# -> r:id="rId6" w:type="default"), plus the harmless re-wrapping of the
# long base64 o:gfxdata blob's embedded newlines. Diffing the two sides
# with attributes canonically sorted shows them to be byte-for-byte
# identical - there is no textual, structural, formatting or style
# change in this commit (matching the generic "Moving from 2.0.1 to
# 2.0.2" version-bump commit message, which is unrelated to this fixture
# file's content).
#
# So the correct edit is a no-op against the Word object model: we must
# not alter any text, run/paragraph formatting, styles, headers,
# footnotes, section properties or drawings.
$d = $word.ActiveDocument
